$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing cell values (recalculated means/stats for updated date ranges) ---
$ws.Range("G2").Value = 0.784300305683267
$ws.Range("G5").Value = 866.217137874454
$ws.Range("G6").Value = 866.217137874454
$ws.Range("G7").Value = 866.217137874454
$ws.Range("G8").Value = 866.217137874454
$ws.Range("G9").Value = 0.0271396316459539
$ws.Range("G10").Value = 0.0271396316459539
$ws.Range("G13").Value = 1.00386101694915
$ws.Range("G14").Value = 1.00386101694915
$ws.Range("G22").Value = 1088.53917177276
$ws.Range("G23").Value = 1088.53917177276
$ws.Range("G24").Value = 1088.53917177276
$ws.Range("G25").Value = 1088.53917177276
$ws.Range("G26").Value = 0.0310605180215409
$ws.Range("G27").Value = 0.0310605180215409
$ws.Range("G30").Value = 0.96791186440678
$ws.Range("G31").Value = 0.96791186440678
$ws.Range("G47").Value = 0.872261016949153
$ws.Range("L47").Value = 0.5696
$ws.Range("G48").Value = 0.872261016949153
$ws.Range("L48").Value = 0.5696
$ws.Range("G53").Value = 1.76624039540715
$ws.Range("H53").Value = 8.683298584136059
$ws.Range("G64").Value = 0.878221666666667
$ws.Range("L64").Value = 0.6136
$ws.Range("M64").Value = 1.25744
$ws.Range("N64").Value = 1.44627
$ws.Range("G65").Value = 0.878221666666667
$ws.Range("L65").Value = 0.6136
$ws.Range("M65").Value = 1.25744
$ws.Range("N65").Value = 1.44627
$ws.Range("G70").Value = 2.46621133145764
$ws.Range("H70").Value = 8.683298584136059
$ws.Range("G77").Value = 0.0190405407865001
$ws.Range("G78").Value = 0.0190405407865001
$ws.Range("G81").Value = 0.920033333333333
$ws.Range("I81").Value = 1.8742
$ws.Range("M81").Value = 1.35114
$ws.Range("G82").Value = 0.920033333333333
$ws.Range("I82").Value = 1.8742
$ws.Range("M82").Value = 1.35114
$ws.Range("G88").Value = 2.3754853871695
$ws.Range("H88").Value = 8.683298584136059
$ws.Range("G96").Value = 0.013071529207761
$ws.Range("L96").Value = 0.00424
$ws.Range("G97").Value = 0.013071529207761
$ws.Range("L97").Value = 0.00424
$ws.Range("F101").Value = 0.9788
$ws.Range("G101").Value = 0.9764
$ws.Range("M101").Value = 1.41756
$ws.Range("N101").Value = 1.75638
$ws.Range("F102").Value = 0.9788
$ws.Range("G102").Value = 0.9764
$ws.Range("M102").Value = 1.41756
$ws.Range("N102").Value = 1.75638
$ws.Range("G108").Value = 2.268332577916
$ws.Range("H108").Value = 8.683298584136059
$ws.Range("G116").Value = 0.0093518486391456
$ws.Range("L116").Value = 0.00409
$ws.Range("G117").Value = 0.0093518486391456
$ws.Range("L117").Value = 0.00409
$ws.Range("F121").Value = 0.9788
$ws.Range("G121").Value = 0.963983333333333
$ws.Range("I121").Value = 1.8742
$ws.Range("M121").Value = 1.41756
$ws.Range("F122").Value = 0.9788
$ws.Range("G122").Value = 0.963983333333333
$ws.Range("I122").Value = 1.8742
$ws.Range("M122").Value = 1.41756
$ws.Range("G128").Value = 1.95637605617687
$ws.Range("H128").Value = 8.683298584136059
$ws.Range("G136").Value = 0.0108523045122107
$ws.Range("L136").Value = 0.00459
$ws.Range("G137").Value = 0.0108523045122107
$ws.Range("L137").Value = 0.00459
$ws.Range("G141").Value = 1.05269333333333
$ws.Range("I141").Value = 1.8742
$ws.Range("G142").Value = 1.05269333333333
$ws.Range("I142").Value = 1.8742
$ws.Range("G156").Value = 0.0117059111375401
$ws.Range("L156").Value = 0.00469
$ws.Range("G157").Value = 0.0117059111375401
$ws.Range("L157").Value = 0.00469
$ws.Range("G161").Value = 1.032295
$ws.Range("G162").Value = 1.032295
$ws.Range("G176").Value = 0.0135205624302073
$ws.Range("G177").Value = 0.0135205624302073
$ws.Range("G196").Value = 0.0135628027289897
$ws.Range("G197").Value = 0.0135628027289897
$ws.Range("G216").Value = 0.014303922551772
$ws.Range("G217").Value = 0.014303922551772

# --- Append new rows 227-245 for the 2019 - 2023 reporting period ---
$ws.Range("A227").Value = "Waikawa at Huritini"
$ws.Range("B227").Value = "ASPM"
$ws.Range("C227").Value = "D"
$ws.Range("D227").Value = "2019 - 2023"
$ws.Range("E227").Value = "RepSite"
$ws.Range("F227").Value = 0.263
$ws.Range("G227").Value = 0.2854
$ws.Range("H227").Value = 0.436
$ws.Range("I227").Value = 0.436
$ws.Range("L227").Value = 0.263
$ws.Range("M227").Value = 0.4024
$ws.Range("N227").Value = 0.436
$ws.Range("O227").Value = 1782982
$ws.Range("P227").Value = 5493786
$ws.Range("Q227").Value = "Horowhenua District"
$ws.Range("R227").Value = "Waiopehu"
$ws.Range("S227").Value = "Waikawa"
$ws.Range("T227").Value = "West_9a"
$ws.Range("A228").Value = "Waikawa at Huritini"
$ws.Range("B228").Value = "DRP (95th Percentile)"
$ws.Range("C228").Value = "B"
$ws.Range("D228").Value = "2019 - 2023"
$ws.Range("E228").Value = "RepSite"
$ws.Range("F228").Value = 0.015
$ws.Range("G228").Value = 0.0159074074074074
$ws.Range("H228").Value = 0.034
$ws.Range("I228").Value = 0.0238
$ws.Range("L228").Value = 0.019
$ws.Range("M228").Value = 0.021
$ws.Range("N228").Value = 0.023
$ws.Range("O228").Value = 1782982
$ws.Range("P228").Value = 5493786
$ws.Range("Q228").Value = "Horowhenua District"
$ws.Range("R228").Value = "Waiopehu"
$ws.Range("S228").Value = "Waikawa"
$ws.Range("T228").Value = "West_9a"
$ws.Range("U228").Value = "mg/L"
$ws.Range("A229").Value = "Waikawa at Huritini"
$ws.Range("B229").Value = "DRP (Median)"
$ws.Range("C229").Value = "C"
$ws.Range("D229").Value = "2019 - 2023"
$ws.Range("E229").Value = "RepSite"
$ws.Range("F229").Value = 0.015
$ws.Range("G229").Value = 0.0159074074074074
$ws.Range("H229").Value = 0.034
$ws.Range("I229").Value = 0.0238
$ws.Range("L229").Value = 0.019
$ws.Range("M229").Value = 0.021
$ws.Range("N229").Value = 0.023
$ws.Range("O229").Value = 1782982
$ws.Range("P229").Value = 5493786
$ws.Range("Q229").Value = "Horowhenua District"
$ws.Range("R229").Value = "Waiopehu"
$ws.Range("S229").Value = "Waikawa"
$ws.Range("T229").Value = "West_9a"
$ws.Range("U229").Value = "mg/L"
$ws.Range("A230").Value = "Waikawa at Huritini"
$ws.Range("B230").Value = "E coli (>260)"
$ws.Range("C230").Value = "E"
$ws.Range("D230").Value = "2019 - 2023"
$ws.Range("E230").Value = "RepSite"
$ws.Range("F230").Value = 371
$ws.Range("G230").Value = 758.185185185185
$ws.Range("H230").Value = 6200
$ws.Range("I230").Value = 2764.4
$ws.Range("J230").Value = 31.4814814814815
$ws.Range("K230").Value = 72.2222222222222
$ws.Range("L230").Value = 640
$ws.Range("M230").Value = 1409.44
$ws.Range("N230").Value = 1996.48
$ws.Range("O230").Value = 1782982
$ws.Range("P230").Value = 5493786
$ws.Range("Q230").Value = "Horowhenua District"
$ws.Range("R230").Value = "Waiopehu"
$ws.Range("S230").Value = "Waikawa"
$ws.Range("T230").Value = "West_9a"
$ws.Range("U230").Value = "% exceedances over 260/100 mL"
$ws.Range("A231").Value = "Waikawa at Huritini"
$ws.Range("B231").Value = "E coli (>540)"
$ws.Range("C231").Value = "E"
$ws.Range("D231").Value = "2019 - 2023"
$ws.Range("E231").Value = "RepSite"
$ws.Range("F231").Value = 371
$ws.Range("G231").Value = 758.185185185185
$ws.Range("H231").Value = 6200
$ws.Range("I231").Value = 2764.4
$ws.Range("J231").Value = 31.4814814814815
$ws.Range("K231").Value = 72.2222222222222
$ws.Range("L231").Value = 640
$ws.Range("M231").Value = 1409.44
$ws.Range("N231").Value = 1996.48
$ws.Range("O231").Value = 1782982
$ws.Range("P231").Value = 5493786
$ws.Range("Q231").Value = "Horowhenua District"
$ws.Range("R231").Value = "Waiopehu"
$ws.Range("S231").Value = "Waikawa"
$ws.Range("T231").Value = "West_9a"
$ws.Range("U231").Value = "% exceedances over 540/100 mL"
$ws.Range("A232").Value = "Waikawa at Huritini"
$ws.Range("B232").Value = "E coli (Median)"
$ws.Range("C232").Value = "E"
$ws.Range("D232").Value = "2019 - 2023"
$ws.Range("E232").Value = "RepSite"
$ws.Range("F232").Value = 371
$ws.Range("G232").Value = 758.185185185185
$ws.Range("H232").Value = 6200
$ws.Range("I232").Value = 2764.4
$ws.Range("J232").Value = 31.4814814814815
$ws.Range("K232").Value = 72.2222222222222
$ws.Range("L232").Value = 640
$ws.Range("M232").Value = 1409.44
$ws.Range("N232").Value = 1996.48
$ws.Range("O232").Value = 1782982
$ws.Range("P232").Value = 5493786
$ws.Range("Q232").Value = "Horowhenua District"
$ws.Range("R232").Value = "Waiopehu"
$ws.Range("S232").Value = "Waikawa"
$ws.Range("T232").Value = "West_9a"
$ws.Range("U232").Value = "E. coli/100 mL"
$ws.Range("A233").Value = "Waikawa at Huritini"
$ws.Range("B233").Value = "E coli (95th Percentile)"
$ws.Range("C233").Value = "E"
$ws.Range("D233").Value = "2019 - 2023"
$ws.Range("E233").Value = "RepSite"
$ws.Range("F233").Value = 371
$ws.Range("G233").Value = 758.185185185185
$ws.Range("H233").Value = 6200
$ws.Range("I233").Value = 2764.4
$ws.Range("J233").Value = 31.4814814814815
$ws.Range("K233").Value = 72.2222222222222
$ws.Range("L233").Value = 640
$ws.Range("M233").Value = 1409.44
$ws.Range("N233").Value = 1996.48
$ws.Range("O233").Value = 1782982
$ws.Range("P233").Value = 5493786
$ws.Range("Q233").Value = "Horowhenua District"
$ws.Range("R233").Value = "Waiopehu"
$ws.Range("S233").Value = "Waikawa"
$ws.Range("T233").Value = "West_9a"
$ws.Range("U233").Value = "E. coli/100 mL"
$ws.Range("A234").Value = "Waikawa at Huritini"
$ws.Range("B234").Value = "MCI"
$ws.Range("C234").Value = "C"
$ws.Range("D234").Value = "2019 - 2023"
$ws.Range("E234").Value = "RepSite"
$ws.Range("F234").Value = 98
$ws.Range("G234").Value = 99.244
$ws.Range("H234").Value = 116
$ws.Range("I234").Value = 116
$ws.Range("L234").Value = 98
$ws.Range("M234").Value = 112.899
$ws.Range("N234").Value = 116
$ws.Range("O234").Value = 1782982
$ws.Range("P234").Value = 5493786
$ws.Range("Q234").Value = "Horowhenua District"
$ws.Range("R234").Value = "Waiopehu"
$ws.Range("S234").Value = "Waikawa"
$ws.Range("T234").Value = "West_9a"
$ws.Range("A235").Value = "Waikawa at Huritini"
$ws.Range("B235").Value = "Ammoniacal-N (95th Percentile)"
$ws.Range("C235").Value = "A"
$ws.Range("D235").Value = "2019 - 2023"
$ws.Range("E235").Value = "RepSite"
$ws.Range("F235").Value = 0.01336
$ws.Range("G235").Value = 0.0138028150817662
$ws.Range("H235").Value = 0.0364551084982659
$ws.Range("I235").Value = 0.02854
$ws.Range("L235").Value = 0.0133
$ws.Range("M235").Value = 0.01972
$ws.Range("N235").Value = 0.02413
$ws.Range("O235").Value = 1782982
$ws.Range("P235").Value = 5493786
$ws.Range("Q235").Value = "Horowhenua District"
$ws.Range("R235").Value = "Waiopehu"
$ws.Range("S235").Value = "Waikawa"
$ws.Range("T235").Value = "West_9a"
$ws.Range("U235").Value = "mg NH4-N/L"
$ws.Range("A236").Value = "Waikawa at Huritini"
$ws.Range("B236").Value = "Ammoniacal-N (Median)"
$ws.Range("C236").Value = "A"
$ws.Range("D236").Value = "2019 - 2023"
$ws.Range("E236").Value = "RepSite"
$ws.Range("F236").Value = 0.01336
$ws.Range("G236").Value = 0.0138028150817662
$ws.Range("H236").Value = 0.0364551084982659
$ws.Range("I236").Value = 0.02854
$ws.Range("L236").Value = 0.0133
$ws.Range("M236").Value = 0.01972
$ws.Range("N236").Value = 0.02413
$ws.Range("O236").Value = 1782982
$ws.Range("P236").Value = 5493786
$ws.Range("Q236").Value = "Horowhenua District"
$ws.Range("R236").Value = "Waiopehu"
$ws.Range("S236").Value = "Waikawa"
$ws.Range("T236").Value = "West_9a"
$ws.Range("U236").Value = "mg NH4-N/L"
$ws.Range("A237").Value = "Waikawa at Huritini"
$ws.Range("B237").Value = "Nitrate-N (95th Percentile)"
$ws.Range("C237").Value = "A"
$ws.Range("D237").Value = "2019 - 2023"
$ws.Range("E237").Value = "RepSite"
$ws.Range("F237").Value = 0.8835
$ws.Range("G237").Value = 0.876259259259259
$ws.Range("H237").Value = 1.54
$ws.Range("I237").Value = 1.43
$ws.Range("L237").Value = 0.832
$ws.Range("M237").Value = 1.2028
$ws.Range("N237").Value = 1.3254
$ws.Range("O237").Value = 1782982
$ws.Range("P237").Value = 5493786
$ws.Range("Q237").Value = "Horowhenua District"
$ws.Range("R237").Value = "Waiopehu"
$ws.Range("S237").Value = "Waikawa"
$ws.Range("T237").Value = "West_9a"
$ws.Range("U237").Value = "mg NO3-N/L"
$ws.Range("A238").Value = "Waikawa at Huritini"
$ws.Range("B238").Value = "Nitrate-N (Median)"
$ws.Range("C238").Value = "A"
$ws.Range("D238").Value = "2019 - 2023"
$ws.Range("E238").Value = "RepSite"
$ws.Range("F238").Value = 0.8835
$ws.Range("G238").Value = 0.876259259259259
$ws.Range("H238").Value = 1.54
$ws.Range("I238").Value = 1.43
$ws.Range("L238").Value = 0.832
$ws.Range("M238").Value = 1.2028
$ws.Range("N238").Value = 1.3254
$ws.Range("O238").Value = 1782982
$ws.Range("P238").Value = 5493786
$ws.Range("Q238").Value = "Horowhenua District"
$ws.Range("R238").Value = "Waiopehu"
$ws.Range("S238").Value = "Waikawa"
$ws.Range("T238").Value = "West_9a"
$ws.Range("U238").Value = "mg NO3-N/L"
$ws.Range("A239").Value = "Waikawa at Huritini"
$ws.Range("B239").Value = "QMCI"
$ws.Range("C239").Value = "C"
$ws.Range("D239").Value = "2019 - 2023"
$ws.Range("E239").Value = "RepSite"
$ws.Range("F239").Value = 5.02
$ws.Range("G239").Value = 4.7164
$ws.Range("H239").Value = 6.262
$ws.Range("I239").Value = 6.262
$ws.Range("L239").Value = 5.02
$ws.Range("M239").Value = 5.8658
$ws.Range("N239").Value = 6.262
$ws.Range("O239").Value = 1782982
$ws.Range("P239").Value = 5493786
$ws.Range("Q239").Value = "Horowhenua District"
$ws.Range("R239").Value = "Waiopehu"
$ws.Range("S239").Value = "Waikawa"
$ws.Range("T239").Value = "West_9a"
$ws.Range("A240").Value = "Waikawa at Huritini"
$ws.Range("B240").Value = "Soluble Inorganic Nitrogen (95th Percentile)"
$ws.Range("D240").Value = "2019 - 2023"
$ws.Range("E240").Value = "RepSite"
$ws.Range("F240").Value = 0.928
$ws.Range("G240").Value = 0.914111111111111
$ws.Range("H240").Value = 1.569
$ws.Range("I240").Value = 1.4734
$ws.Range("L240").Value = 0.87
$ws.Range("M240").Value = 1.25432
$ws.Range("N240").Value = 1.40126
$ws.Range("O240").Value = 1782982
$ws.Range("P240").Value = 5493786
$ws.Range("Q240").Value = "Horowhenua District"
$ws.Range("R240").Value = "Waiopehu"
$ws.Range("S240").Value = "Waikawa"
$ws.Range("T240").Value = "West_9a"
$ws.Range("U240").Value = "g/m3"
$ws.Range("A241").Value = "Waikawa at Huritini"
$ws.Range("B241").Value = "Soluble Inorganic Nitrogen (Median)"
$ws.Range("D241").Value = "2019 - 2023"
$ws.Range("E241").Value = "RepSite"
$ws.Range("F241").Value = 0.928
$ws.Range("G241").Value = 0.914111111111111
$ws.Range("H241").Value = 1.569
$ws.Range("I241").Value = 1.4734
$ws.Range("L241").Value = 0.87
$ws.Range("M241").Value = 1.25432
$ws.Range("N241").Value = 1.40126
$ws.Range("O241").Value = 1782982
$ws.Range("P241").Value = 5493786
$ws.Range("Q241").Value = "Horowhenua District"
$ws.Range("R241").Value = "Waiopehu"
$ws.Range("S241").Value = "Waikawa"
$ws.Range("T241").Value = "West_9a"
$ws.Range("U241").Value = "g/m3"
$ws.Range("A242").Value = "Waikawa at Huritini"
$ws.Range("B242").Value = "Total Nitrogen (95th Percentile)"
$ws.Range("D242").Value = "2019 - 2023"
$ws.Range("E242").Value = "RepSite"
$ws.Range("F242").Value = 1.07
$ws.Range("G242").Value = 1.0537037037037
$ws.Range("H242").Value = 1.65
$ws.Range("I242").Value = 1.624
$ws.Range("L242").Value = 1.04
$ws.Range("M242").Value = 1.3796
$ws.Range("N242").Value = 1.4524
$ws.Range("O242").Value = 1782982
$ws.Range("P242").Value = 5493786
$ws.Range("Q242").Value = "Horowhenua District"
$ws.Range("R242").Value = "Waiopehu"
$ws.Range("S242").Value = "Waikawa"
$ws.Range("T242").Value = "West_9a"
$ws.Range("U242").Value = "g/m3"
$ws.Range("A243").Value = "Waikawa at Huritini"
$ws.Range("B243").Value = "Total Nitrogen (Median)"
$ws.Range("D243").Value = "2019 - 2023"
$ws.Range("E243").Value = "RepSite"
$ws.Range("F243").Value = 1.07
$ws.Range("G243").Value = 1.0537037037037
$ws.Range("H243").Value = 1.65
$ws.Range("I243").Value = 1.624
$ws.Range("L243").Value = 1.04
$ws.Range("M243").Value = 1.3796
$ws.Range("N243").Value = 1.4524
$ws.Range("O243").Value = 1782982
$ws.Range("P243").Value = 5493786
$ws.Range("Q243").Value = "Horowhenua District"
$ws.Range("R243").Value = "Waiopehu"
$ws.Range("S243").Value = "Waikawa"
$ws.Range("T243").Value = "West_9a"
$ws.Range("U243").Value = "g/m3"
$ws.Range("A244").Value = "Waikawa at Huritini"
$ws.Range("B244").Value = "Total Phosphorus (95th Percentile)"
$ws.Range("D244").Value = "2019 - 2023"
$ws.Range("E244").Value = "RepSite"
$ws.Range("F244").Value = 0.036
$ws.Range("G244").Value = 0.0400185185185185
$ws.Range("H244").Value = 0.099
$ws.Range("I244").Value = 0.0646
$ws.Range("L244").Value = 0.042
$ws.Range("M244").Value = 0.05132
$ws.Range("N244").Value = 0.059
$ws.Range("O244").Value = 1782982
$ws.Range("P244").Value = 5493786
$ws.Range("Q244").Value = "Horowhenua District"
$ws.Range("R244").Value = "Waiopehu"
$ws.Range("S244").Value = "Waikawa"
$ws.Range("T244").Value = "West_9a"
$ws.Range("U244").Value = "g/m3"
$ws.Range("A245").Value = "Waikawa at Huritini"
$ws.Range("B245").Value = "Total Phosphorus (Median)"
$ws.Range("D245").Value = "2019 - 2023"
$ws.Range("E245").Value = "RepSite"
$ws.Range("F245").Value = 0.036
$ws.Range("G245").Value = 0.0400185185185185
$ws.Range("H245").Value = 0.099
$ws.Range("I245").Value = 0.0646
$ws.Range("L245").Value = 0.042
$ws.Range("M245").Value = 0.05132
$ws.Range("N245").Value = 0.059
$ws.Range("O245").Value = 1782982
$ws.Range("P245").Value = 5493786
$ws.Range("Q245").Value = "Horowhenua District"
$ws.Range("R245").Value = "Waiopehu"
$ws.Range("S245").Value = "Waikawa"
$ws.Range("T245").Value = "West_9a"
$ws.Range("U245").Value = "g/m3"
